# Correcao pvs placa 04A
#
# The "Hxy" channel label in column B had three inconsistent spellings
# scattered through the sheet ("A_Hxyy" / "A_Hxyx" / "A_Hxy" and their
# "C_" counterparts). This normalises every occurrence to the single
# correct spelling "A_Hxy" / "C_Hxy".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "A_Hxy"
$ws.Range("B22").Value = "C_Hxy"
$ws.Range("B26").Value = "A_Hxy"
$ws.Range("B27").Value = "C_Hxy"
$ws.Range("B30").Value = "A_Hxy"
$ws.Range("B31").Value = "C_Hxy"
$ws.Range("B38").Value = "A_Hxy"
$ws.Range("B39").Value = "C_Hxy"
$ws.Range("B44").Value = "A_Hxy"
$ws.Range("B45").Value = "A_Hxy"
$ws.Range("B47").Value = "C_Hxy"
$ws.Range("B48").Value = "C_Hxy"
$ws.Range("B54").Value = "A_Hxy"
$ws.Range("B55").Value = "C_Hxy"
$ws.Range("B62").Value = "A_Hxy"
$ws.Range("B63").Value = "C_Hxy"
$ws.Range("B66").Value = "A_Hxy"
$ws.Range("B67").Value = "C_Hxy"
$ws.Range("B71").Value = "A_Hxy"
$ws.Range("B72").Value = "C_Hxy"

# Leave the sheet scrolled to the top with C12 selected, matching the
# saved view state of the edited workbook.
$ws.Range("C12").Select()
